# Updated cryptos list on Fri Mar 22 02:24:27 UTC 2024 with GitHub Actions
#
# Refreshes the Price (col D) and Volume(1h) (col E) values scraped from
# coinranking.com, and corrects three rows where the scraped ranking order
# shifted position (Uniswap/WrappedBTC, Toncoin/Litecoin, PEPE/Maker), which
# also moves the Coin name (col B) and Link (col C) for those rows.
#
# Price cells are forced back to text (NumberFormat "@") before the write so
# Excel doesn't silently coerce numeric-looking strings (e.g. "564.32",
# "0.0000266") into floating point doubles and mangle their literal
# representation; the style is reset to Normal afterwards so no new cell
# style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.324.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.462.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.017.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.457.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.280.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.983"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "408.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "29.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "611.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -14.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.355.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.97%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0780"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.63%  "
$ws.Range("E40").Value = "  -6.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0409"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
